$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.420.19'
$ws.Range('E2').Value = '  +3.23%  '
$ws.Range('D3').Value = '1.589.47'
$ws.Range('E3').Value = '  +1.63%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.91%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.58'
$ws.Range('E5').Value = '  +1.32%  '
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('E7').Value = '  +1.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.09'
$ws.Range('E8').Value = '  +6.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.250'
$ws.Range('E9').Value = '  +0.24%  '
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0884'
$ws.Range('E11').Value = '  +1.92%  '
$ws.Range('D12').Value = '1.815.17'
$ws.Range('E12').Value = '  +1.56%  '
$ws.Range('D13').Value = '1.584.24'
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('E14').Value = '  +1.97%  '
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').Value = '28.381.92'
$ws.Range('E16').Value = '  +3.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.24'
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.58'
$ws.Range('E18').Value = '  +2.24%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.49'
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0706'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  +0.98%  '
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.75'
$ws.Range('E25').Value = '  +1.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.21'
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('E28').Value = '  -1.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.13'
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = '1.397.28'
$ws.Range('E34').Value = '  -3.94%  '
$ws.Range('E35').Value = '  -1.38%  '
$ws.Range('E36').Value = '  -8.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.34'
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('E39').Value = '  +9.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.540'
$ws.Range('E40').Value = '  -0.33%  '
$ws.Range('E42').Value = '  +1.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.88'
$ws.Range('E43').Value = '  +2.43%  '
$ws.Range('E44').Value = '  -2.40%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.978'
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.22'
$ws.Range('D47').Value = '1.725.44'
$ws.Range('E47').Value = '  +1.43%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.15'
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('B49').Value = 'mCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.13'
$ws.Range('E49').Value = '  +1.25%  '
$ws.Range('D50').Value = '0.0₆0103'
$ws.Range('E50').Value = '  +16.28%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '41.51'
$ws.Range('E51').Value = '  +9.96%  '
